$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)

# The header paragraph originally reads (as a sequence of runs):
#   "吾立方外包管理网页客户端软件 v" [eastAsia] + "2" [default] + <bookmark _GoBack>
#   + "021" [default] + "." [eastAsia] + "0.0" [default] + <tab> + <PAGE field>
#
# It needs to become:
#   "吾立方" [eastAsia] + "Nuke" [default] + "插件 v2021.4.1" [eastAsia]
#   + <tab> + <PAGE field>
#
# and the "_GoBack" bookmark needs to move out of the header and into the
# (otherwise empty) first body paragraph.
#
# Edits are applied back-to-front (highest offsets first) so earlier
# character offsets stay valid while later ones shift.

# 1) "." + "0.0" (chars 20-24) -> "插件 v2021.4.1" (keeps the eastAsia run's
#    formatting, since that run anchors the replacement range).
$r1 = $hdr.Range.Duplicate
$r1.SetRange(20, 24)
$r1.Text = "插件 v2021.4.1"

# 2) "2" + bookmark + "021" (chars 16-20) -> "Nuke". This span straddles the
#    _GoBack bookmark, which removes it from the header along with the text.
$r2 = $hdr.Range.Duplicate
$r2.SetRange(16, 20)
$r2.Text = "Nuke"

# 3) Leading run (chars 0-16) -> "吾立方".
$r3 = $hdr.Range.Duplicate
$r3.SetRange(0, 16)
$r3.Text = "吾立方"

# 4) Re-create the "_GoBack" bookmark around the (empty) first body paragraph.
$d.Bookmarks.Add("_GoBack", $d.Paragraphs.Item(1).Range)
